$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cached/source data values ---
# Columns: B = NewReno, C = Reno, D = Tahoe, E = Vegas (rows 2-11 correspond to CBR 1-10)
$ws.Range("B2").Value = 412.683570592
$ws.Range("B3").Value = 412.686103742
$ws.Range("B4").Value = 412.703686247
$ws.Range("B5").Value = 412.597377572
$ws.Range("B6").Value = 412.296111748
$ws.Range("B7").Value = 412.162141265
$ws.Range("B8").Value = 411.05910602
$ws.Range("B9").Value = 391.117122205
$ws.Range("B10").Value = 396.894532528
$ws.Range("B11").Value = 383.308538579

$ws.Range("C2").Value = 412.683570592
$ws.Range("C3").Value = 412.686103742
$ws.Range("C4").Value = 412.703686247
$ws.Range("C5").Value = 412.597377572
$ws.Range("C6").Value = 412.296111748
$ws.Range("C7").Value = 412.102087991
$ws.Range("C8").Value = 410.761151765
$ws.Range("C9").Value = 395.423732841
$ws.Range("C10").Value = 395.494557688
$ws.Range("C11").Value = 383.589059373

$ws.Range("D2").Value = 412.683570592
$ws.Range("D3").Value = 412.686103742
$ws.Range("D4").Value = 412.703686247
$ws.Range("D5").Value = 412.597377572
$ws.Range("D6").Value = 412.296111748
$ws.Range("D7").Value = 412.278257449
$ws.Range("D8").Value = 410.961142562
$ws.Range("D9").Value = 390.274996917
$ws.Range("D10").Value = 394.596011269
$ws.Range("D11").Value = 382.571450063

$ws.Range("E2").Value = 398.371256884
$ws.Range("E3").Value = 397.707204364
$ws.Range("E4").Value = 397.712537974
$ws.Range("E5").Value = 397.690984017
$ws.Range("E6").Value = 397.590142968
$ws.Range("E7").Value = 396.936397683
$ws.Range("E8").Value = 395.339424357
$ws.Range("E9").Value = 364.100167794
$ws.Range("E10").Value = 363.423255107
$ws.Range("E11").Value = 452.743370103

# --- Update chart value axis minimum scale ---
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.MinimumScale = 360

# --- Update active selection ---
[void]$ws.Range("K19").Select()
